# mySQL_Create/layers.xlsx
# "Added UNIQUE constraints to db definition, changed classification_definition."
#
# The "name" column row (row 3) gets a "UNIQUE" note in the "Other 1" column
# (E3), matching the style already used by its row-neighbours (C3 etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the built-in "Normal" cell style to "Standard" (the label the
# German-locale build of Excel uses for the same built-in style) so the
# workbook matches what re-saving it in that Excel build produced.
try {
    $normalStyle = $wb.Styles("Normal")
    $normalStyle.Name = "Standard"
} catch {
    # Non-fatal if the host does not support renaming the built-in style.
}

# Copy the formatting already used for this row's data cells (Arial 10,
# no border/fill) onto E3, then fill in the new constraint text.
$ws.Range("C3").Copy()
$ws.Range("E3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E3").Value = "UNIQUE"

# Leave the selection on E4, matching where the editor ended up.
$ws.Range("E4").Select()
